$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1457.7358
$ws.Range("J17").Value = 1457.7358
$ws.Range("L17").Value = 4373.207399999999
$ws.Range("N17").Value = -4709.207399999999

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H43").Value = 1098.6666
$ws.Range("J43").Value = 1098.6666
$ws.Range("L43").Value = 1098.6666
$ws.Range("N43").Value = -1236.6666

$ws.Range("H62").Value = 922.5
$ws.Range("I62").Value = 896.6667
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 896.6667
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = -272.6667
$ws.Range("N62").Value = -2248

$ws.Range("H65").Value = 922.5
$ws.Range("I65").Value = 896.6667
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 4483.3335
$ws.Range("L65").Value = 5000
$ws.Range("M65").Value = -1363.3335
$ws.Range("N65").Value = -11240

$ws.Range("H106").Value = 3786.5
$ws.Range("I106").Value = 3983.125
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 3983.125
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -3352.125
$ws.Range("N106").Value = -4262

$ws.Range("H137").Value = 243684.61
$ws.Range("I137").Value = 1850.5714
$ws.Range("J137").Value = 582252.25
$ws.Range("K137").Value = 5551.7142
$ws.Range("L137").Value = 1746756.75
$ws.Range("M137").Value = -3001.7142
$ws.Range("N137").Value = -1751856.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 742.3333
$ws.Range("I2").Value = 666.75
$ws.Range("J2").Value = 893.5
$ws.Range("K2").Value = 666.75
$ws.Range("L2").Value = 893.5
$ws.Range("M2").Value = -553.75
$ws.Range("N2").Value = -1119.5

$ws.Range("H32").Value = 7183.4106
$ws.Range("I32").Value = 3904.9
$ws.Range("K32").Value = 3904.9
$ws.Range("M32").Value = -3617.9

$ws.Range("H61").Value = 36855.035
$ws.Range("I61").Value = 2159.318
$ws.Range("K61").Value = 2159.318
$ws.Range("M61").Value = -1947.318

$ws.Range("H74").Value = 46036.086
$ws.Range("I74").Value = 72970.36
$ws.Range("K74").Value = 72970.36
$ws.Range("M74").Value = -72096.36

$ws.Range("H77").Value = 46036.086
$ws.Range("I77").Value = 72970.36
$ws.Range("K77").Value = 364851.8
$ws.Range("M77").Value = -360483.8

$ws.Range("H116").Value = 742.3333
$ws.Range("I116").Value = 666.75
$ws.Range("J116").Value = 893.5
$ws.Range("K116").Value = 666.75
$ws.Range("L116").Value = 893.5
$ws.Range("M116").Value = 1627.25
$ws.Range("N116").Value = -5481.5

$ws.Range("H132").Value = 2673.5557
$ws.Range("I132").Value = 2252
$ws.Range("K132").Value = 6756
$ws.Range("M132").Value = -4226

$ws.Range("H136").Value = 36855.035
$ws.Range("I136").Value = 2159.318
$ws.Range("K136").Value = 6477.954000000001
$ws.Range("M136").Value = -3927.954000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 742.3333
$ws.Range("I3").Value = 666.75
$ws.Range("J3").Value = 893.5
$ws.Range("K3").Value = 666.75
$ws.Range("L3").Value = 893.5
$ws.Range("M3").Value = -552.75
$ws.Range("N3").Value = -1121.5

$ws.Range("H132").Value = 45366.57
$ws.Range("J132").Value = 45366.57
$ws.Range("L132").Value = 45366.57
$ws.Range("N132").Value = -55486.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 90294.5
$ws.Range("J68").Value = 90294.5
$ws.Range("L68").Value = 90294.5
$ws.Range("N68").Value = -91792.5

$ws.Range("H71").Value = 90294.5
$ws.Range("J71").Value = 90294.5
$ws.Range("L71").Value = 270883.5
$ws.Range("N71").Value = -278371.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 92
$ws.Range("I2").Value = 146.2
$ws.Range("K2").Value = 877.1999999999999
$ws.Range("M2").Value = -764.1999999999999

$ws.Range("H97").Value = 124.75
$ws.Range("J97").Value = 99.5
$ws.Range("L97").Value = 298.5
$ws.Range("N97").Value = -1290.5

$ws.Range("H121").Value = 3345.8
$ws.Range("I121").Value = 1615
$ws.Range("J121").Value = 4499.6665
$ws.Range("K121").Value = 4845
$ws.Range("L121").Value = 13498.9995
$ws.Range("M121").Value = -3535
$ws.Range("N121").Value = -16118.9995

$ws.Range("H131").Value = 1300.25
$ws.Range("I131").Value = 1083.3
$ws.Range("K131").Value = 3249.9
$ws.Range("M131").Value = 1790.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3897197
$ws.Range("I11").Value = 1737992.2
$ws.Range("J11").Value = 5336667
$ws.Range("K11").Value = 1737992.2
$ws.Range("L11").Value = 5336667
$ws.Range("M11").Value = -1737853.2
$ws.Range("N11").Value = -5336945

$ws.Range("H21").Value = 732635.6
$ws.Range("I21").Value = 2502974.8
$ws.Range("K21").Value = 2502974.8
$ws.Range("M21").Value = -2502801.8

$ws.Range("H24").Value = 21023.629
$ws.Range("I24").Value = 8800
$ws.Range("J24").Value = 21493.77
$ws.Range("K24").Value = 8800
$ws.Range("L24").Value = 21493.77
$ws.Range("M24").Value = -8627
$ws.Range("N24").Value = -21839.77

$ws.Range("H30").Value = 732635.6
$ws.Range("I30").Value = 2502974.8
$ws.Range("K30").Value = 2502974.8
$ws.Range("M30").Value = -2502869.8

$ws.Range("H33").Value = 1006515.2
$ws.Range("J33").Value = 7239.1113
$ws.Range("L33").Value = 7239.1113
$ws.Range("N33").Value = -7743.1113

$ws.Range("H97").Value = 1650.24
$ws.Range("I97").Value = 613.6
$ws.Range("J97").Value = 5796.8
$ws.Range("K97").Value = 613.6
$ws.Range("L97").Value = 5796.8
$ws.Range("M97").Value = -117.6
$ws.Range("N97").Value = -6788.8

$ws.Range("H126").Value = 6755.222
$ws.Range("I126").Value = 3599.8
$ws.Range("K126").Value = 10799.4
$ws.Range("M126").Value = -8329.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4277471
$ws.Range("I40").Value = 4190.2
$ws.Range("K40").Value = 4190.2
$ws.Range("M40").Value = -4054.2

$ws.Range("H55").Value = 6956.25
$ws.Range("I55").Value = 1427.8334
$ws.Range("J55").Value = 23541.5
$ws.Range("K55").Value = 1427.8334
$ws.Range("L55").Value = 23541.5
$ws.Range("M55").Value = -1254.8334
$ws.Range("N55").Value = -23887.5

$ws.Range("H68").Value = 3187.2
$ws.Range("I68").Value = 3487.125
$ws.Range("J68").Value = 1987.5
$ws.Range("K68").Value = 3487.125
$ws.Range("L68").Value = 1987.5
$ws.Range("M68").Value = -2738.125
$ws.Range("N68").Value = -3485.5

$ws.Range("H71").Value = 3187.2
$ws.Range("I71").Value = 3487.125
$ws.Range("J71").Value = 1987.5
$ws.Range("K71").Value = 17435.625
$ws.Range("L71").Value = 9937.5
$ws.Range("M71").Value = -13691.625
$ws.Range("N71").Value = -17425.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("K81").Value = 2000
$ws.Range("M81").Value = -939

$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("K84").Value = 10000
$ws.Range("M84").Value = -4696
